$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9:85 down to 10:86
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(9, 3).Value = "Maule"
$ws.Cells.Item(9, 4).Value = 44550
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = 100112022
$ws.Cells.Item(9, 7).Value = "Arveja Verde"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 300
$ws.Cells.Item(9, 11).Value = 16000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 16000
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Carahue"
$ws.Cells.Item(9, 16).Value = 640
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
